$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest cryptos data refresh

$ws.Range("D2").Value = "26.953.08"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "1.558.72"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("E8").Value = "  +0.27%  "

$ws.Range("E9").Value = "  +0.50%  "

$ws.Range("E10").Value = "  +1.65%  "

$ws.Range("E11").Value = "  -0.47%  "

$ws.Range("D12").Value = "1.780.53"
$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("D13").Value = "1.555.38"
$ws.Range("E13").Value = "  -0.10%  "

$ws.Range("E14").Value = "  +0.03%  "

$ws.Range("E15").Value = "  +0.08%  "

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.91"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("D17").Value = "26.956.47"
$ws.Range("E17").Value = "  +0.02%  "

$ws.Range("E18").Value = "  +1.71%  "

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.11"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  -0.77%  "

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.39"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +1.04%  "

$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("E22").Value = "  +1.50%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.23"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +0.44%  "

$ws.Range("E24").Value = "  -1.20%  "

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.67"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -0.68%  "

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.61"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -0.59%  "

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.09"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +1.00%  "

$ws.Range("E28").Value = "  +1.14%  "

$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("E30").Value = "  +0.78%  "

$ws.Range("E31").Value = "  +2.72%  "

$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("E33").Value = "  +2.53%  "

$ws.Range("D34").Value = "1.427.12"
$ws.Range("E34").Value = "  +0.31%  "

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.60"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  +1.11%  "

$ws.Range("E36").Value = "  +8.21%  "

$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.32"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  +1.83%  "

$ws.Range("E38").Value = "  +1.05%  "

$ws.Range("E39").Value = "  +2.23%  "

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.86"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +2.75%  "

$ws.Range("E41").Value = "  -0.47%  "

$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("E43").Value = "  +1.76%  "

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.31"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +0.72%  "

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.60"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -0.16%  "

$ws.Range("E46").Value = "  +0.52%  "

$ws.Range("D47").Value = "1.694.26"
$ws.Range("E47").Value = "  +0.03%  "

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.34"
$ws.Range("D48").Style = $origStyle

$ws.Range("D49").Value = "0.0₆0104"
$ws.Range("E49").Value = "  +4.09%  "

$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("E51").Value = "  +0.12%  "
